# Auto-generated edit script: refreshes Leve profit-calculation values
# (currentAveragePrice / LevePrice / LeveProfit columns) across the
# crafting-job worksheets, per the scheduled market-data runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 11682.875
$ws.Range("I62").Value = 16092.143
$ws.Range("J62").Value = 5509.9
$ws.Range("K62").Value = 16092.143
$ws.Range("L62").Value = 5509.9
$ws.Range("M62").Value = -15468.143
$ws.Range("N62").Value = -6757.9

# Row 65
$ws.Range("H65").Value = 11682.875
$ws.Range("I65").Value = 16092.143
$ws.Range("J65").Value = 5509.9
$ws.Range("K65").Value = 80460.715
$ws.Range("L65").Value = 27549.5
$ws.Range("M65").Value = -77340.715
$ws.Range("N65").Value = -33789.5

# Row 116
$ws.Range("H116").Value = 3750.75
$ws.Range("I116").Value = 4555.8
$ws.Range("J116").Value = 2821.8462
$ws.Range("K116").Value = 4555.8
$ws.Range("L116").Value = 2821.8462
$ws.Range("M116").Value = -1113.8
$ws.Range("N116").Value = -9705.8462

# Row 137
$ws.Range("H137").Value = 2127.3936
$ws.Range("I137").Value = 2083.186
$ws.Range("J137").Value = 2233
$ws.Range("K137").Value = 6249.558000000001
$ws.Range("L137").Value = 6699
$ws.Range("M137").Value = -3699.558000000001
$ws.Range("N137").Value = -11799

$ws = $wb.Worksheets.Item("ARM")
# Row 36
$ws.Range("H36").Value = 5931.5
$ws.Range("I36").Value = 1863
$ws.Range("K36").Value = 1863
$ws.Range("M36").Value = -1517

# Row 61
$ws.Range("H61").Value = 297981.44
$ws.Range("I61").Value = 221358.81
$ws.Range("J61").Value = 480416.28
$ws.Range("K61").Value = 221358.81
$ws.Range("L61").Value = 480416.28
$ws.Range("M61").Value = -221146.81
$ws.Range("N61").Value = -480840.28

# Row 63
$ws.Range("H63").Value = 3345.1
$ws.Range("I63").Value = 3277.9443
$ws.Range("J63").Value = 3949.5
$ws.Range("K63").Value = 3277.9443
$ws.Range("L63").Value = 3949.5
$ws.Range("M63").Value = -2591.9443
$ws.Range("N63").Value = -5321.5

# Row 66
$ws.Range("H66").Value = 3345.1
$ws.Range("I66").Value = 3277.9443
$ws.Range("J66").Value = 3949.5
$ws.Range("K66").Value = 16389.7215
$ws.Range("L66").Value = 19747.5
$ws.Range("M66").Value = -12957.7215
$ws.Range("N66").Value = -26611.5

# Row 74
$ws.Range("H74").Value = 112428.73
$ws.Range("I74").Value = 118327.18
$ws.Range("J74").Value = 76616.71000000001
$ws.Range("K74").Value = 118327.18
$ws.Range("L74").Value = 76616.71000000001
$ws.Range("M74").Value = -117453.18
$ws.Range("N74").Value = -78364.71000000001

# Row 77
$ws.Range("H77").Value = 112428.73
$ws.Range("I77").Value = 118327.18
$ws.Range("J77").Value = 76616.71000000001
$ws.Range("K77").Value = 591635.8999999999
$ws.Range("L77").Value = 383083.55
$ws.Range("M77").Value = -587267.8999999999
$ws.Range("N77").Value = -391819.55

# Row 102
$ws.Range("H102").Value = 2717.625
$ws.Range("I102").Value = 1656
$ws.Range("J102").Value = 4487
$ws.Range("K102").Value = 1656
$ws.Range("L102").Value = 4487
$ws.Range("M102").Value = -34
$ws.Range("N102").Value = -7731

# Row 110
$ws.Range("H110").Value = 1279.9
$ws.Range("I110").Value = 1316.5
$ws.Range("J110").Value = 1225
$ws.Range("K110").Value = 1316.5
$ws.Range("L110").Value = 1225
$ws.Range("M110").Value = 728.5
$ws.Range("N110").Value = -5315

# Row 122
$ws.Range("H122").Value = 4670.5947
$ws.Range("I122").Value = 4645.5483
$ws.Range("K122").Value = 13936.6449
$ws.Range("M122").Value = -11486.6449

# Row 132
$ws.Range("H132").Value = 3025.9805
$ws.Range("I132").Value = 2791.3784
$ws.Range("J132").Value = 3646
$ws.Range("K132").Value = 8374.135200000001
$ws.Range("L132").Value = 10938
$ws.Range("M132").Value = -5844.135200000001
$ws.Range("N132").Value = -15998

# Row 136
$ws.Range("H136").Value = 297981.44
$ws.Range("I136").Value = 221358.81
$ws.Range("J136").Value = 480416.28
$ws.Range("K136").Value = 664076.4299999999
$ws.Range("L136").Value = 1441248.84
$ws.Range("M136").Value = -661526.4299999999
$ws.Range("N136").Value = -1446348.84

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1414.5
$ws.Range("I20").Value = 1183.1364
$ws.Range("J20").Value = 1838.6666
$ws.Range("K20").Value = 1183.1364
$ws.Range("L20").Value = 1838.6666
$ws.Range("M20").Value = -936.1364000000001
$ws.Range("N20").Value = -2332.6666

# Row 94
$ws.Range("H94").Value = 3099.6667
$ws.Range("I94").Value = 1533
$ws.Range("J94").Value = 4666.3335
$ws.Range("K94").Value = 1533
$ws.Range("L94").Value = 4666.3335
$ws.Range("M94").Value = -1082
$ws.Range("N94").Value = -5568.3335

# Row 105
$ws.Range("H105").Value = 2079.3704
$ws.Range("I105").Value = 1535.5
$ws.Range("K105").Value = 1535.5
$ws.Range("M105").Value = 211.5

# Row 134
$ws.Range("H134").Value = 2507.5493
$ws.Range("I134").Value = 1912.2307
$ws.Range("J134").Value = 4136.8423
$ws.Range("K134").Value = 5736.6921
$ws.Range("L134").Value = 12410.5269
$ws.Range("M134").Value = -3201.6921
$ws.Range("N134").Value = -17480.5269

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 95.36
$ws.Range("I7").Value = 114
$ws.Range("J7").Value = 75.166664
$ws.Range("K7").Value = 114
$ws.Range("L7").Value = 75.166664
$ws.Range("M7").Value = -1
$ws.Range("N7").Value = -301.166664

# Row 122
$ws.Range("H122").Value = 1772.2142
$ws.Range("I122").Value = 1053
$ws.Range("J122").Value = 2059.9
$ws.Range("K122").Value = 3159
$ws.Range("L122").Value = 6179.700000000001
$ws.Range("M122").Value = -709
$ws.Range("N122").Value = -11079.7

# Row 134
$ws.Range("H134").Value = 1978.8286
$ws.Range("I134").Value = 1253.25
$ws.Range("J134").Value = 2946.2666
$ws.Range("K134").Value = 3759.75
$ws.Range("L134").Value = 8838.799800000001
$ws.Range("M134").Value = -1224.75
$ws.Range("N134").Value = -13908.7998

$ws = $wb.Worksheets.Item("CUL")
# Row 132
$ws.Range("H132").Value = 3747.3333
$ws.Range("I132").Value = 627
$ws.Range("J132").Value = 5619.533
$ws.Range("K132").Value = 5643
$ws.Range("L132").Value = 50575.79700000001
$ws.Range("M132").Value = -3113
$ws.Range("N132").Value = -55635.79700000001

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 1528.5
$ws.Range("I122").Value = 1528.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4585.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2135.5
$ws.Range("N122").ClearContents()  # HQ profit column no longer populated for this listing

$ws = $wb.Worksheets.Item("LTW")
# Row 26
$ws.Range("H26").Value = 6054.5
$ws.Range("I26").Value = 5109
$ws.Range("J26").Value = 7000
$ws.Range("K26").Value = 5109
$ws.Range("L26").Value = 7000
$ws.Range("M26").Value = -4814
$ws.Range("N26").Value = -7590

# Row 46
$ws.Range("H46").Value = 988.7692
$ws.Range("I46").Value = 983.3333
$ws.Range("J46").Value = 1001
$ws.Range("K46").Value = 983.3333
$ws.Range("L46").Value = 1001
$ws.Range("M46").Value = -795.3333
$ws.Range("N46").Value = -1377

# Row 68
$ws.Range("H68").Value = 3230
$ws.Range("I68").Value = 2965
$ws.Range("J68").Value = 3318.3333
$ws.Range("K68").Value = 2965
$ws.Range("L68").Value = 3318.3333
$ws.Range("M68").Value = -2216
$ws.Range("N68").Value = -4816.3333

# Row 71
$ws.Range("H71").Value = 3230
$ws.Range("I71").Value = 2965
$ws.Range("J71").Value = 3318.3333
$ws.Range("K71").Value = 14825
$ws.Range("L71").Value = 16591.6665
$ws.Range("M71").Value = -11081
$ws.Range("N71").Value = -24079.6665

# Row 93
$ws.Range("H93").Value = 1429
$ws.Range("I93").Value = 1357.5
$ws.Range("J93").Value = 1595.8334
$ws.Range("K93").Value = 1357.5
$ws.Range("L93").Value = 1595.8334
$ws.Range("M93").Value = -109.5
$ws.Range("N93").Value = -4091.8334

# Row 122
$ws.Range("H122").Value = 2880
$ws.Range("I122").Value = 3450
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 10350
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -7900
$ws.Range("N122").Value = -12400

# Row 132
$ws.Range("H132").Value = 13515.056
$ws.Range("I132").Value = 4813.8335
$ws.Range("J132").Value = 17865.666
$ws.Range("K132").Value = 14441.5005
$ws.Range("L132").Value = 53596.99800000001
$ws.Range("M132").Value = -11911.5005
$ws.Range("N132").Value = -58656.99800000001

$ws = $wb.Worksheets.Item("WVR")
# Row 32
$ws.Range("H32").Value = 7575
$ws.Range("I32").Value = 3862.5
$ws.Range("J32").Value = 15000
$ws.Range("K32").Value = 3862.5
$ws.Range("L32").Value = 15000
$ws.Range("M32").Value = -3545.5
$ws.Range("N32").Value = -15634

# Row 33
$ws.Range("H33").Value = 5634.25
$ws.Range("I33").Value = 3719
$ws.Range("J33").Value = 5907.857
$ws.Range("K33").Value = 3719
$ws.Range("L33").Value = 5907.857
$ws.Range("M33").Value = -3469
$ws.Range("N33").Value = -6407.857

# Row 36
$ws.Range("H36").Value = 5634.25
$ws.Range("I36").Value = 3719
$ws.Range("J36").Value = 5907.857
$ws.Range("K36").Value = 3719
$ws.Range("L36").Value = 5907.857
$ws.Range("M36").Value = -3469
$ws.Range("N36").Value = -6407.857

# Row 132
$ws.Range("H132").Value = 1979.4359
$ws.Range("I132").Value = 1351.0435
$ws.Range("J132").Value = 2882.75
$ws.Range("K132").Value = 4053.1305
$ws.Range("L132").Value = 8648.25
$ws.Range("M132").Value = -1523.1305
$ws.Range("N132").Value = -13708.25
